$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldText = "Administrator, Miss Dina Nasr"
$newText = "Miss Dina Nasr, Administrator"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value -eq $oldText) {
        $cell.Value = $newText
    }
}
